$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test")

# Select the Test named range on the Test sheet (matches final cursor state there)
$ws.Range("D5:G10").Select()

# Insert the new "Conversion" worksheet right after "Test" and make it active
$cv = $wb.Worksheets.Add($null, $ws)
$cv.Name = "Conversion"

# --- Header row 6 : bold/shaded header style copied from Test!D5:G5 ---
$cv.Range("B6").Value = "AAA"
$cv.Range("C6").Value = "BBB"
$cv.Range("D6").Value = "CCC"
$cv.Range("E6").Value = "DDD"
$ws.Range("D5:G5").Copy()
$cv.Range("B6:E6").PasteSpecial(-4122)

# --- Pre-format date/time cells so they pick up the same number formats
#     already used on the Test sheet (custom dd/mm/yyyy hh:mm:ss, and the
#     built-in m/d/yy h:mm date-time format) ---
$cv.Range("E7").NumberFormat = "m/d/yy h:mm"
$cv.Range("C11").NumberFormat = "m/d/yy h:mm"
$cv.Range("E8").NumberFormat = "dd/mm/yyyy\ hh:mm:ss"
$cv.Range("E9").NumberFormat = "dd/mm/yyyy\ hh:mm:ss"
$cv.Range("E10").NumberFormat = "dd/mm/yyyy\ hh:mm:ss"

# --- Row 7 ---
$cv.Range("B7").Value = "'-14.65"
$cv.Range("C7").Value = "hello"
$cv.Range("D7").Value = $true
$ws.Range("F6").Copy()
$cv.Range("D7").PasteSpecial(-4122)
$cv.Range("E7").Value = 30692.5

# --- Row 8 ---
$cv.Range("B8").Value = "not-a-number"
$cv.Range("C8").Value = 42.24
$cv.Range("D8").Value = "'TRUE"
$cv.Range("E8").Value = "'06.02.2012 16:15:23"

# --- Row 9 ---
$cv.Range("B9").Value = "'11.7"
$cv.Range("C9").Value = $true
$cv.Range("E9").Value = "'11.01.1984 12:00:00"

# --- Row 10 ---
$cv.Range("B10").Value = 780.9
$cv.Range("D10").Value = "not-a-boolean"
$cv.Range("E10").Value = "not-a-date"

# --- Row 11 ---
$cv.Range("C11").Value = 30692.5
$cv.Range("D11").Value = "'FALSE"
$cv.Range("E11").Value = 357.67

# --- Column widths (bestFit approximations, quantised to the engine's
#     nearest 1/6-character column-width grid) ---
$cv.Columns("B").ColumnWidth = 18.307291666666668
$cv.Columns("C").ColumnWidth = 14.307291666666666
$cv.Columns("D").ColumnWidth = 14.307291666666666
$cv.Columns("E").ColumnWidth = 17.022135416666668

# --- Page setup: portrait, 2cm top/bottom margins (matches Test sheet) ---
$cv.PageSetup.Orientation = 1
$cv.PageSetup.TopMargin = 56.692913399999995
$cv.PageSetup.BottomMargin = 56.692913399999995

# --- Workbook-level defined name for the new region ---
$wb.Names.Add("Conversion", "=Conversion!`$B`$6:`$E`$11")

# --- Final selection on the now-active Conversion sheet ---
$cv.Range("A3").Select()
